$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3571.5454
$ws.Range("I80").Value = 3431.5
$ws.Range("K80").Value = 10294.5
$ws.Range("M80").Value = -9296.5
$ws.Range("H83").Value = 3571.5454
$ws.Range("I83").Value = 3431.5
$ws.Range("K83").Value = 30883.5
$ws.Range("M83").Value = -25891.5
$ws.Range("H100").Value = 2163.2856
$ws.Range("I100").Value = 3136.5
$ws.Range("K100").Value = 3136.5
$ws.Range("M100").Value = -2595.5
$ws.Range("H118").Value = 395
$ws.Range("I118").Value = 395
$ws.Range("K118").Value = 1185
$ws.Range("M118").Value = 472

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8550.392
$ws.Range("I32").Value = 8935.772000000001
$ws.Range("K32").Value = 8935.772000000001
$ws.Range("M32").Value = -8648.772000000001
$ws.Range("H45").Value = 2625.3333
$ws.Range("J45").Value = 1070
$ws.Range("L45").Value = 1070
$ws.Range("N45").Value = -1824
$ws.Range("H74").Value = 5000.364
$ws.Range("I74").Value = 4500.5
$ws.Range("J74").Value = 9999
$ws.Range("K74").Value = 4500.5
$ws.Range("L74").Value = 9999
$ws.Range("M74").Value = -3626.5
$ws.Range("N74").Value = -11747
$ws.Range("H77").Value = 5000.364
$ws.Range("I77").Value = 4500.5
$ws.Range("J77").Value = 9999
$ws.Range("K77").Value = 22502.5
$ws.Range("L77").Value = 49995
$ws.Range("M77").Value = -18134.5
$ws.Range("N77").Value = -58731
$ws.Range("H122").Value = 4299.25
$ws.Range("I122").Value = 3871.818
$ws.Range("K122").Value = 11615.454
$ws.Range("M122").Value = -9165.454000000002
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1333.3334
$ws.Range("J11").Value = 1750
$ws.Range("L11").Value = 1750
$ws.Range("N11").Value = -2030
$ws.Range("H99").Value = 2199
$ws.Range("I99").Value = 2199
$ws.Range("K99").Value = 2199
$ws.Range("M99").Value = -701

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2126.9092
$ws.Range("I31").Value = 1466.5
$ws.Range("K31").Value = 1466.5
$ws.Range("M31").Value = -1171.5
$ws.Range("H34").Value = 2126.9092
$ws.Range("I34").Value = 1466.5
$ws.Range("K34").Value = 1466.5
$ws.Range("M34").Value = -1264.5
$ws.Range("H105").Value = 4541.857
$ws.Range("I105").Value = 5633.3335
$ws.Range("J105").Value = 3723.25
$ws.Range("K105").Value = 5633.3335
$ws.Range("L105").Value = 3723.25
$ws.Range("M105").Value = -3886.3335
$ws.Range("N105").Value = -7217.25
$ws.Range("H107").Value = 1226.591
$ws.Range("I107").Value = 1478.4667
$ws.Range("J107").Value = 686.8570999999999
$ws.Range("K107").Value = 1478.4667
$ws.Range("L107").Value = 686.8570999999999
$ws.Range("M107").Value = 441.5333000000001
$ws.Range("N107").Value = -4526.8571

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("H81").Value = 5747.25
$ws.Range("J81").Value = 7495
$ws.Range("L81").Value = 22485
$ws.Range("N81").Value = -24731
$ws.Range("H84").Value = 5747.25
$ws.Range("J84").Value = 7495
$ws.Range("L84").Value = 67455
$ws.Range("N84").Value = -78687
$ws.Range("H98").Value = 1997.8889
$ws.Range("J98").Value = 1229.6
$ws.Range("L98").Value = 3688.8
$ws.Range("N98").Value = -6684.799999999999
$ws.Range("H129").Value = 1968.5555
$ws.Range("I129").Value = 1290
$ws.Range("J129").Value = 2053.375
$ws.Range("K129").Value = 3870
$ws.Range("L129").Value = 6160.125
$ws.Range("M129").Value = 1130
$ws.Range("N129").Value = -16160.125
$ws.Range("N46").ClearContents()

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 350
$ws.Range("I36").Value = 200
$ws.Range("K36").Value = 200
$ws.Range("M36").Value = 285
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H102").Value = 2990.3635
$ws.Range("I102").Value = 2801.375
$ws.Range("K102").Value = 2801.375
$ws.Range("M102").Value = -1179.375
$ws.Range("H122").Value = 4470
$ws.Range("I122").Value = 4470
$ws.Range("K122").Value = 13410
$ws.Range("M122").Value = -10960
$ws.Range("H126").Value = 5621.8335
$ws.Range("I126").Value = 5963.4
$ws.Range("K126").Value = 17890.2
$ws.Range("M126").Value = -15420.2
$ws.Range("N58").ClearContents()

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 9999.5
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10226
$ws.Range("H28").Value = 9999.5
$ws.Range("J28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("N28").Value = -10464
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("H35").Value = 1100
$ws.Range("I35").Value = 1100
$ws.Range("K35").Value = 1100
$ws.Range("M35").Value = -764
$ws.Range("H37").Value = 9999.5
$ws.Range("J37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("N37").Value = -10214
$ws.Range("H40").Value = 4333.3335
$ws.Range("I40").Value = 4333.3335
$ws.Range("K40").Value = 4333.3335
$ws.Range("M40").Value = -4197.3335
$ws.Range("H42").Value = 49464
$ws.Range("J42").Value = 49464
$ws.Range("L42").Value = 49464
$ws.Range("N42").Value = -50590
$ws.Range("H46").Value = 6555.5557
$ws.Range("I46").Value = 2000
$ws.Range("K46").Value = 2000
$ws.Range("M46").Value = -1812
$ws.Range("H49").Value = 49464
$ws.Range("J49").Value = 49464
$ws.Range("L49").Value = 49464
$ws.Range("N49").Value = -49758
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 800
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 800
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -598
$ws.Range("N61").Value = -1604
$ws.Range("H100").Value = 3120
$ws.Range("I100").Value = 3120
$ws.Range("K100").Value = 3120
$ws.Range("M100").Value = -2579
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -5540
$ws.Range("N33").ClearContents()

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("H122").Value = 6185
$ws.Range("I122").Value = 6185
$ws.Range("K122").Value = 18555
$ws.Range("M122").Value = -16105
$ws.Range("H132").Value = 3303.2856
$ws.Range("I132").Value = 2031.75
$ws.Range("K132").Value = 6095.25
$ws.Range("M132").Value = -3565.25
$ws.Range("N19").ClearContents()
